# Updated cryptos list values per the upstream GitHub Actions scrape commit.
# Column D ("Price") is always stored as text in this sheet, so each write
# forces text via NumberFormat "@" (else Excel auto-coerces numeric-looking
# strings like "1.004" or "1.000" into numbers and drops formatting/precision),
# then clears the format again so no stray style id is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.728.35"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3884"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("E9").Value = "  -4.44%  "
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07139"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -7.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.055"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.741.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.924"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001046"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06583"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "78.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.70"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.139"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.764.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.386"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.98"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.72"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.273"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.940.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.264"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -11.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.059"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.713"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08693"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.89"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.15%  "
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.511"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02252"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.62%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.053"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06044"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6389"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2079"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.186"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.812"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.806"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5889"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.961"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.56%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06922"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.28%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.138"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.34%  "
